# Append the new PR log entry (#45) as row 21 of the PR log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21
$ws.Cells.Item($row, 1).Value = 45
$ws.Cells.Item($row, 2).Value = "trying squash"
$ws.Cells.Item($row, 3).Value = "riya-morankar"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "edit1 to main"
# Leading apostrophe forces the date-looking string to stay text (matching
# the other Date-column cells), rather than being auto-converted to a date
# serial number.
$ws.Cells.Item($row, 6).Value = "'2025-06-18"
